$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.055.76"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.551"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.288"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0714"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.048.43"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.25"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.797.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.624"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.058.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.80"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.13"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.113"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0515"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.67"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.51"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.398.69"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.655"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0188"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.22"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.921"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.78"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.15"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.87%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.07"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0137"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +10.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0496"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.947.81"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.07%  "
